$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lake Champlain (Vermont) code fixed from "CL" to "LC"
$ws.Range("C30").Value = "LC"

# Columbus (Ohio) code fixed from "CB" to "CO"
$ws.Range("C34").Value = "CO"

# Update the remembered selection to match the saved workbook state
$ws.Range("F29").Select()
